# Finals GT vs RR
# Enter match results for the Finals row (row 86) on Sheet1, enter the
# "Rank 1" tally for the Finals round (column H, rows 96-102) on Sheet1,
# and record the winner predictions (Rank 1 / Rank 2) for the Finals on
# Sheet2. Then fix up the view state (active sheet / selections) and the
# display format picked up on U102 as a side effect of the edits.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: Finals (row 86) raw scores for each of the 7 players ---
$ws1.Range("E86").Value = 40
$ws1.Range("H86").Value = 80
$ws1.Range("K86").Value = 60
$ws1.Range("N86").Value = 100
$ws1.Range("Q86").Value = 70
$ws1.Range("T86").Value = 50
$ws1.Range("W86").Value = 0

# --- Sheet1: Finals "Rank 1" guess counts for the scorecard (column H) ---
$ws1.Range("H96").Value = 0
$ws1.Range("H97").Value = 3
$ws1.Range("H98").Value = 3
$ws1.Range("H99").Value = 11
$ws1.Range("H100").Value = 0
$ws1.Range("H101").Value = 0
$ws1.Range("H102").Value = 0

# --- Sheet2: Finals GT vs RR winner predictions (Rank 1 / Rank 2) ---
$ws2.Range("J40").Value = "Sibi"
$ws2.Range("J41").Value = "Justin"

# U102 picked up extra decimal precision (manually bumped via Increase
# Decimal in the UI) - replicate the resulting custom number format /
# style so the exported styles.xml matches.
$ws1.Range("U102").NumberFormat = "0.000000"

# --- View state: Sheet2 selection moves off J41, and Sheet1 becomes the
# active/selected tab again. ---
$ws2.Activate()
$ws2.Range("F47:G47").Select()
$ws1.Activate()
$ws1.Range("U103").Select()
